# Self Service price sheet enhancement
# Updates PO Cost / Matrix Cost values on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: PO Cost (G3) drops to 7.5
$ws.Range("G3").Value = 7.5

# Matrix Cost (I3:I5) increases to 12.5
$ws.Range("I3").Value = 12.5
$ws.Range("I4").Value = 12.5
$ws.Range("I5").Value = 12.5

# Leave the selection on the last edited cell, matching the saved view state.
$null = $ws.Range("I5").Select()
